$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.232.43"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.445.29"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +4.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.31"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.46%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.573"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.43%  "

$ws.Range("E9").Value = "  +5.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.84"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.08%  "

$ws.Range("E11").Value = "  +2.52%  "

$ws.Range("E12").Value = "  -1.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.00"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.870.54"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.090.80"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.65%  "

$ws.Range("E16").Value = "  +5.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.434.35"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.18%  "

$ws.Range("E18").Value = "  +7.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.13"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.94"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.19%  "

$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.77"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.02%  "

$ws.Range("E24").Value = "  +2.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.60"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0798"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +8.27%  "

$ws.Range("E29").Value = "  +4.13%  "

$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.33"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.04%  "

$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.14"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.84"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.17%  "

$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.33"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.29"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.46%  "

$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("E38").Value = "  +0.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "40.17"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.418"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +10.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "319.21"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +8.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.75"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "143.57"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0967"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0528"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.49"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.84%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.576"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.23%  "

$ws.Range("B48").Value = "Polygon"
$ws.Range("C48").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.408"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +7.06%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0227"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.04"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("E51").Value = "  +4.96%  "
